$wb = $excel.ActiveWorkbook
$wsControl = $wb.Worksheets.Item("control condition")
$wsModel = $wb.Worksheets.Item("model-supported condition")

# control condition
$wsControl.Range("K2").Value = 2
$wsControl.Range("Q2").Value = 2
$wsControl.Range("Q3").Value = 2
$wsControl.Range("K4").Value = 1
$wsControl.Range("K5").Value = 2
$wsControl.Range("N5").Value = 1
$wsControl.Range("Q5").Value = 2
$wsControl.Range("K6").Value = 2
$wsControl.Range("N6").Value = 1
$wsControl.Range("Q7").Value = 1
$wsControl.Range("K8").Value = 1
$wsControl.Range("N8").Value = 1
$wsControl.Range("Q8").Value = 1
$wsControl.Range("N10").Value = 1
$wsControl.Range("Q10").Value = 2
$wsControl.Range("Q11").Value = 1
$wsControl.Range("Q15").Value = 2
$wsControl.Range("K16").Value = 2
$wsControl.Range("N16").Value = 1
$wsControl.Range("Q16").Value = 2
$wsControl.Range("K17").Value = 1
$wsControl.Range("N17").Value = 2
$wsControl.Range("Q17").Value = 2
$wsControl.Range("K20").Value = 2
$wsControl.Range("Q20").Value = 2
$wsControl.Range("Q25").Value = 2

# model-supported condition
$wsModel.Range("N2").Value = 1
$wsModel.Range("Q2").Value = 2
$wsModel.Range("K3").Value = 2
$wsModel.Range("Q3").Value = 2
$wsModel.Range("K4").Value = 2
$wsModel.Range("Q4").Value = 2
$wsModel.Range("K6").Value = 2
$wsModel.Range("Q6").Value = 2
$wsModel.Range("K8").Value = 2
$wsModel.Range("Q8").Value = 2
$wsModel.Range("Q11").Value = 1
$wsModel.Range("Q14").Value = 2
$wsModel.Range("N15").Value = 1
$wsModel.Range("K16").Value = 2
$wsModel.Range("Q16").Value = 2
$wsModel.Range("K17").Value = 1
$wsModel.Range("Q17").Value = 2
$wsModel.Range("K19").Value = 2
$wsModel.Range("Q19").Value = 2
$wsModel.Range("K20").Value = 1
$wsModel.Range("Q21").Value = 1

# Update sheet views / selections to match target state.
# model-supported condition: selection becomes M2:R25, and it is no longer the active tab.
$wsModel.Activate() | Out-Null
$wsModel.Range("M2:R25").Select() | Out-Null

# control condition: becomes the active tab, selection becomes J2:R25.
$wsControl.Activate() | Out-Null
$wsControl.Range("J2:R25").Select() | Out-Null

